# DSW2 update - incomplete. encaptulation still fucking
$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("DSW1")
$ws5 = $wb.Worksheets.Item("DSW2")

# --- DSW2 (sheet5): fill in the addressing table that was started ---

# Header row
$ws5.Range("C1").Value = "IPv4"
$ws5.Range("D1").Value = "SubnetMask"
$ws5.Range("E1").Value = "IPv6"

# Gi1/0/1
$ws5.Range("A2").Value = "Gi1/0/1"
$ws5.Range("B2").Value = "VLAN5"
$ws5.Range("C2").Value = "10.1.4.10"
$ws5.Range("D2").Value = 255255255252
$ws5.Range("E2").Value = "2026:4::2/64"
$ws5.Range("G2").Value = "Tilføj IPv6"

# Gi1/0/23
$ws5.Range("A3").Value = "Gi1/0/23"
$ws5.Range("B3").Value = "VLAN1"
$ws5.Range("C3").Value = "10.2.4.14"
$ws5.Range("D3").Value = 255255255252
$ws5.Range("E3").Value = "2026:3::2/64"
$ws5.Range("G3").Value = "Har ""ekstra"" IPv6 2026:3::1"

# Loopback2
$ws5.Range("A4").Value = "Loopback2"
$ws5.Range("C4").Value = "22.22.22.22"
$ws5.Range("D4").Value = 255255255252
$ws5.Range("E4").Value = "FEC0:22::2/64"
$ws5.Range("G4").Value = "Tilføj IPv6?"

# VLAN1
$ws5.Range("A5").Value = "VLAN1"
$ws5.Range("C5").Value = "10.2.4.14"
$ws5.Range("D5").Value = 255255255252
$ws5.Range("E5").Value = "2026:3::2/64"
$ws5.Range("G5").Value = "Har ""ekstra"" IPv6 2026:3::1"

# VLAN5
$ws5.Range("A6").Value = "VLAN5"
$ws5.Range("C6").Value = "10.1.4.10"
$ws5.Range("D6").Value = 255255255252
$ws5.Range("E6").Value = "2026:4::2/64"

# VLAN10
$ws5.Range("A7").Value = "VLAN10"
$ws5.Range("C7").Value = "10.2.1.2"
$ws5.Range("D7").Value = "255.255.255.0"
$ws5.Range("E7").Value = "2026:10::2"

# VLAN20
$ws5.Range("A8").Value = "VLAN20"
$ws5.Range("C8").Value = "10.2.2.2"
$ws5.Range("D8").Value = "255.255.255.0"
$ws5.Range("E8").Value = "2026:20::2"

# VLAN150
$ws5.Range("A9").Value = "VLAN150"
$ws5.Range("C9").Value = "10.2.3.2"
$ws5.Range("D9").Value = 255255255128
$ws5.Range("E9").Value = "No IPv6"

# VLAN200
$ws5.Range("A10").Value = "VLAN200"
$ws5.Range("C10").Value = "192.168.1.130"
$ws5.Range("D10").Value = 255255255224
$ws5.Range("E10").Value = "No IPv6"

# VLAN250
$ws5.Range("A11").Value = "VLAN250"
$ws5.Range("C11").Value = "10.2.3.130"
$ws5.Range("D11").Value = 255255255128
$ws5.Range("E11").Value = "No IPv6"

# Gi1/0/19 trunk
$ws5.Range("A12").Value = "Gi1/0/19"
$ws5.Range("B12").Value = "Trunk"
$ws5.Range("C12").Value = "10,20,150,200,250"

# Gi1/0/21 trunk
$ws5.Range("A13").Value = "Gi1/0/21"
$ws5.Range("B13").Value = "Trunk"
$ws5.Range("C13").Value = "10,20,150,200,251"

# Subnet mask column uses the thousands-separator number format
$ws5.Range("D2:D11").NumberFormat = "#,##0"

# Column widths (best-fit to content)
$ws5.Columns.Item(1).AutoFit()
$ws5.Columns.Item(3).AutoFit()
$ws5.Columns.Item(4).AutoFit()
$ws5.Columns.Item(5).AutoFit()

# Page setup to match DSW1
$ws5.PageSetup.PaperSize = 9
$ws5.PageSetup.Orientation = 1

# --- DSW1 (sheet4): row-height cleanup ---
$ws4.Rows.Item(4).AutoFit()
$ws4.Rows.Item(5).RowHeight = 15.75

# --- Move the active tab from DSW1 to DSW2 ---
$ws5.Activate()
$null = $ws5.Range("E2").Select()
